$p = $ppt.ActivePresentation
$m = $p.SlideMaster
$nm = $p.NotesMaster

$mTcs = $m.Theme.ThemeColorScheme
$nmTcs = $nm.Theme.ThemeColorScheme

for ($i = 1; $i -le 12; $i++) {
    $mColor = $mTcs.Colors($i)
    $nmColor = $nmTcs.Colors($i)
    $tmp = $mColor.RGB
    $mColor.RGB = $nmColor.RGB
    $nmColor.RGB = $tmp
}
